$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the cryptos list refresh.
# Column D (Price) values that look numeric must be forced to remain
# plain text (matching the original inlineStr cells) instead of being
# auto-coerced into numbers by the COM Value setter, so each such cell
# is briefly switched to a text number format, written, then the
# temporary formatting is cleared again to restore the original style.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.182.01'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +4.10%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.651.45'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +7.42%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.82'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.09'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.645.90'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +7.47%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.607'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.89%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  +3.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.604'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '49.86'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.58%  '
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '693.02'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.233.71'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +7.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.98'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +3.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.705.97'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +8.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '72.283.65'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +4.14%  '
$ws.Range('E19').Value = '  +2.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.50'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.60'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.29%  '
$ws.Range('E22').Value = '  +2.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.85'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +8.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.90'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '103.75'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.03'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.16%  '
$ws.Range('E27').Value = '  +4.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.93'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '35.27'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +3.19%  '
$ws.Range('E30').Value = '  +3.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.39'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +5.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.12'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +15.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '584.44'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.30'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.110'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.65'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.26%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.649.37'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.50%  '
$ws.Range('E39').Value = '  +2.01%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '35.74'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0771'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +6.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.43'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.92%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0463'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +8.78%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.77'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.349'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.84'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +5.66%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.133'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.05%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.45'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.76%  '
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.78'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.02'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +14.67%  '
